$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.7443820224719101
$ws1.Range("C2").Value = 0.9142857142857143
$ws1.Range("D2").Value = 0.5393258426966292
$ws1.Range("E2").Value = 0.6784452296819788
$ws1.Range("F2").Value = 0.587515299877601
$ws1.Range("G2").Value = 0.5479692645444566
$ws1.Range("H2").Value = 0.7443820224719101
$ws1.Range("I2").Value = 288
$ws1.Range("J2").Value = 27
$ws1.Range("K2").Value = 507
$ws1.Range("L2").Value = 246

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.6733067729083665
$ws2.Range("C2").Value = 0.949438202247191
$ws2.Range("D2").Value = 0.7878787878787878

$ws2.Range("B3").Value = 0.9142857142857143
$ws2.Range("C3").Value = 0.5393258426966292
$ws2.Range("D3").Value = 0.6784452296819788

$ws2.Range("B4").Value = 0.7443820224719101
$ws2.Range("C4").Value = 0.7443820224719101
$ws2.Range("D4").Value = 0.7443820224719101
$ws2.Range("E4").Value = 0.7443820224719101

$ws2.Range("B5").Value = 0.7937962435970404
$ws2.Range("C5").Value = 0.7443820224719101
$ws2.Range("D5").Value = 0.7331620087803834

$ws2.Range("B6").Value = 0.7937962435970404
$ws2.Range("C6").Value = 0.7443820224719101
$ws2.Range("D6").Value = 0.7331620087803833

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 507
$ws3.Range("C2").Value = 27
$ws3.Range("B3").Value = 246
$ws3.Range("C3").Value = 288
